# Apply financial data updates to the QADA worksheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("QADA")

# --- Rows where the oldest period (column J, FY2012) became unavailable ---
$ws.Range("J21").Value = "NA"
$ws.Range("J83").Value = "NA"
$ws.Range("J94").Value = "NA"
$ws.Range("J100").Value = "NA"
$ws.Range("J101").Value = "NA"

# --- Revised FY2017 (column E) figures ---
$ws.Range("E27").Value = -15500
$ws.Range("E33").Value = -15500
$ws.Range("E35").Value = -15500
$ws.Range("E81").Value = -15500

# --- Revised "Capital Expenditures" row (row 91) across all periods ---
$ws.Range("D91").Value = -3700
$ws.Range("E91").Value = -3300
$ws.Range("F91").Value = -3200
$ws.Range("G91").Value = -4600
$ws.Range("H91").Value = -4800
$ws.Range("I91").Value = -3100
$ws.Range("J91").Value = -3800
